$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines_B2 = @(
    'SELECT DISTINCT',
    '     std.dbgap_accession AS "dbGaP Accession",',
    '     std.study_name AS "Study Name"',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '    std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female'''
)
$ws.Range("B2").Value = [string]::Join("`r`n", $lines_B2)

$lines_C2 = @(
    'SELECT ',
    '    COUNT(DISTINCT dgn.diagnosis) AS Diagnoses,',
    '    COUNT(DISTINCT prt.participant_id) AS Participants,',
    '    COUNT(DISTINCT std.study_id) AS Studies',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '   std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female'';'
)
$ws.Range("C2").Value = [string]::Join("`r`n", $lines_C2)

$lines_B3 = @(
    'SELECT DISTINCT',
    '    prt.participant_id AS "Participant Id",',
    '    prt.race AS "Race",',
    '    prt.sex_at_birth AS "Sex at Birth",',
    '    std.dbgap_accession AS "dbGaP Accession"',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '    std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female''',
    'ORDER BY ',
    '    prt.participant_id ASC',
    'LIMIT 100;'
)
$ws.Range("B3").Value = [string]::Join("`r`n", $lines_B3)

$lines_B4 = @(
    'SELECT DISTINCT',
    '    prt.participant_id AS "Participant Id",',
    '    dgn.diagnosis_id AS "Diagnosis Id",',
    '    dgn.diagnosis AS "Diagnosis",',
    '    dgn.diagnosis_classification_system AS "Diagnosis Classification System",',
    '    dgn.diagnosis_basis AS "Diagnosis Basis",',
    '    dgn.tumor_classification AS "Tumor Classification",',
    '    dgn.anatomic_site AS "Anatomic Site",',
    '    CASE ',
    '    WHEN dgn.age_at_diagnosis = -999 THEN ''Not Reported''',
    '    WHEN dgn.age_at_diagnosis >= 1000 THEN ',
    '        substr(dgn.age_at_diagnosis, 1, length(dgn.age_at_diagnosis) - 3) || '','' || substr(dgn.age_at_diagnosis, -3)',
    '    ELSE ',
    '        dgn.age_at_diagnosis ',
    'END AS "Age at Diagnosis (days)",',
    '    std.dbgap_accession AS "dbGaP Accession"',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '    std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female''',
    'ORDER BY ',
    '    dgn.diagnosis_id ASC',
    'LIMIT 100;'
)
$ws.Range("B4").Value = [string]::Join("`r`n", $lines_B4)

$lines_B5 = @(
    'SELECT',
    '    DISTINCT prt.participant_id AS "Participant Id",',
    '    trt.treatment_id AS "Treatment Id",',
    '    CASE ',
    '    WHEN trt.age_at_treatment_start = -999 THEN ''Not Reported''',
    '    WHEN trt.age_at_treatment_start >= 1000 THEN ',
    '        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || '','' || substr(trt.age_at_treatment_start, -3)',
    '    ELSE ',
    '        trt.age_at_treatment_start ',
    'END AS "Age at Treatment Start",',
    '    CASE ',
    '    WHEN trt.age_at_treatment_end = -999 THEN ''Not Reported''',
    '    WHEN trt.age_at_treatment_end >= 1000 THEN ',
    '        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || '','' || substr(trt.age_at_treatment_end, -3)',
    '    ELSE ',
    '        trt.age_at_treatment_end ',
    'END AS "Age at Treatment End",',
    '    trt.treatment_type AS "Treatment Type",',
    '    REPLACE(trt.treatment_agent, '';'', '', '') AS "Treatment Agent",',
    '    std.dbgap_accession AS "dbGaP Accession"',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '    std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female''',
    'ORDER BY ',
    '    trt.treatment_id ASC',
    'LIMIT 100;'
)
$ws.Range("B5").Value = [string]::Join("`r`n", $lines_B5)

$lines_B6 = @(
    'SELECT DISTINCT',
    '    prt.participant_id AS "Participant Id",',
    '    trr.treatment_response_id AS "Treatment Response Id",',
    '    trr.response AS "Response",',
    '    CASE ',
    '        WHEN trr.age_at_response = -999 THEN ''Not Reported''',
    '        WHEN trr.age_at_response >= 1000 THEN ',
    '            substr(trr.age_at_response, 1, length(trr.age_at_response) - 3) || '','' || substr(trr.age_at_response, -3)',
    '        ELSE ',
    '            trr.age_at_response ',
    '    END AS "Age at Response",',
    '    trr.response_category AS "Response Category",',
    '    trr.response_system AS "Response System",',
    '    std.dbgap_accession AS "dbGaP Accession"',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '    std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female''',
    'ORDER BY ',
    '    prt.participant_id ASC',
    'LIMIT 100;'
)
$ws.Range("B6").Value = [string]::Join("`r`n", $lines_B6)

$lines_B7 = @(
    'SELECT DISTINCT',
    '    prt.participant_id AS "Participant Id",',
    '    srv.survival_id AS "Survival Id",',
    '    srv.last_known_survival_status AS "Last Known Survival Status",',
    '    CASE ',
    '    WHEN srv.age_at_last_known_survival_status = -999 THEN ''Not Reported''',
    '    WHEN srv.age_at_last_known_survival_status >= 1000 THEN ',
    '        substr(srv.age_at_last_known_survival_status, 1, length(srv.age_at_last_known_survival_status) - 3) || '','' || substr(srv.age_at_last_known_survival_status, -3)',
    '    ELSE ',
    '        srv.age_at_last_known_survival_status ',
    'END AS "Age at Last Known Survival Status",',
    '    srv.first_event AS "First Event",',
    '    srv.cause_of_death AS "Cause of Death",',
    '    std.dbgap_accession AS "dbGaP Accession"',
    'FROM ',
    '    df_study std',
    'LEFT JOIN ',
    '    df_participant prt ON std.study_id = prt."study.study_id"',
    'LEFT JOIN ',
    '    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatments trt ON prt.participant_id = trt."participant.participant_id"',
    'LEFT JOIN ',
    '    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"',
    'LEFT JOIN ',
    '    df_survival srv ON prt.participant_id = srv."participant.participant_id"',
    'LEFT JOIN ',
    '    df_reference_files rfs ON std.study_id = rfs."study.study_id"',
    'WHERE ',
    '    std.dbgap_accession = ''phs003111'' AND prt.sex_at_birth = ''Female''',
    'ORDER BY ',
    '    prt.participant_id ASC',
    'LIMIT 100;'
)
$ws.Range("B7").Value = [string]::Join("`r`n", $lines_B7)
